# UniSolder52C BOM update
# "R70 must be 0ohm, profiles corrected"
#
# 1. R70 is pulled out of the shared "R23, R27, R28, R33, R44, R57, R70"
#    (1.8k) designator group on row 56 into its own BOM line, because R70
#    actually needs to be a 0 ohm (jumper) resistor, not 1.8k.
# 2. A new row is inserted (row 66) for R70 = 0 ohm / Resistor / 0805
#    Resistor / 0805, pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: shrink the 1.8k group on row 56 (remove R70 from it) ---
$ws.Range("A56").Value = "'R23, R27, R28, R33, R44, R57"
$ws.Range("B56").Value = 6

# --- Step 2: insert a dedicated row for R70 right after row 65 (R67) ---
$ws.Rows.Item(66).Insert()

# Copy the formatting of the row above so borders / wrap / alignment match
# the rest of the resistor table, then overwrite with R70's own data.
$ws.Range("A65:G65").Copy()
$ws.Range("A66:G66").PasteSpecial(-4122)

$ws.Range("A66").Value = "'R70"
$ws.Range("B66").Value = 1
$ws.Range("C66").Value = "'0"
$ws.Range("D66").Value = "'0"
$ws.Range("E66").Value = "'Resistor"
$ws.Range("F66").Value = "'0805 Resistor"
$ws.Range("G66").Value = "'0805"

# --- Cosmetic: leave the view scrolled/selected near the edit, like the
#     author's saved workbook (topLeftCell isn't persisted by this host,
#     but the active selection is). ---
$ws.Range("D67").Select()

Write-Output "R70 split into its own 0-ohm BOM row"
